# Updated DEU model - 2025-09-01 09:18
#
# The "wind" fuel row in the fuels sheet is split into two separate fuels:
# "windon" (onshore) and "windoff" (offshore). This is implemented by
# inserting a new row right after the existing "wind" row, renaming the
# original row's commodity to "windon" and filling the new row with
# "windoff" (keeping the same unit as the original "wind" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fuels")

# Insert a new row below the current "wind" row (row 14), pushing every
# following row down by one.
$ws.Rows.Item(15).Insert()

# Row 14 was "wind" / "twh" -> becomes "windon" / "twh"
$ws.Range("C14").Value = "windon"

# New row 15 becomes "windoff" / "twh"
$ws.Range("C15").Value = "windoff"
$ws.Range("E15").Value = "twh"

# Match the cursor position recorded in the saved workbook.
$ws.Activate()
$ws.Range("E16").Select()
